# final commit by Shilpi
# Add a new "E" column of data to the PostMorbidity / PostMorbidityExist
# sheets, and touch the (otherwise still-empty) E column on
# PostMorbidityMissing so its used range grows to match.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PostMorbidity")
$ws1.Range("E2").Value = "1HY_1TC"
$ws1.Range("E3").Value = "2HY_2TC"
$ws1.Range("E4").Value = "3HY_3TC"
$ws1.Range("E5").Value = "4HY_4TC"
$ws1.Range("E6").Value = "5HY_5TC"
$ws1.Range("E7").Value = "6HY_6TC"
$ws1.Range("E8").Value = "7HY_7TC"
$ws1.Range("E9").Value = "8HY_8TC"
$ws1.Range("E10").Value = "9HY_9TC"
$ws1.Range("E11").Value = "10H_10T"

$ws2 = $wb.Worksheets.Item("PostMorbidityExist")
$ws2.Range("E2").Value = "1VA_1VB"
$ws2.Range("E3").Value = "2VA_2VB"
$ws2.Range("E4").Value = "3VA_3VB"
$ws2.Range("E5").Value = "4VA_4VB"
$ws2.Range("E6").Value = "5VA_5VB"

$ws3 = $wb.Worksheets.Item("PostMorbidityMissing")
$ws3.Range("E2:E6").Borders.LineStyle = -4142
$ws3.Range("E2:E6").ClearFormats()
